# toString for DataModel is added
# The VALIDATE() calls used to pass the full words "Error"/"Warning" as the
# severity argument; they now pass the short codes "E"/"W" instead. All the
# cells that stored/echoed that severity (and the concatenated summary
# strings in column G) are updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: literal "Error" -> "E" inside the VALIDATE formula, and the
#     resulting summary text in G2.
$ws.Range("A2").Formula = '=VALIDATE(1=2, "Message A2", C2, "E")'
$ws.Range("G2").Value = "Message A2, E, C2, 1.0, 1.0"

# --- Row 3: literal "Warning" -> "W" inside the VALIDATE formula. G3 keeps
#     the same displayed text ("-") but is rewritten so it mirrors the sheet.
$ws.Range("A3").Formula = '=VALIDATE(1<>2, B3, C3, "W")'
$ws.Range("G3").Value = "-"

# --- Row 4: D4 holds the severity fed into VALIDATE(2=2, B4, C4, D4).
$ws.Range("D4").Value = "E"
$ws.Range("G4").Value = "-"

# --- Row 5: D5 holds the severity fed into VALIDATE(E5=F5, B5, C5, D5).
$ws.Range("D5").Value = "W"
$ws.Range("G5").Value = "Message A5, W, C5, true, true"

# --- Row 6: D6 holds the severity fed into VALIDATE(E6, B6, C6, D6).
$ws.Range("D6").Value = "E"
$ws.Range("G6").Value = "Message A6, E, C6, =IF(0=10,TRUE,FALSE), false"

# --- Row 7: D7 holds the severity fed into VALIDATE(TRUE, B7, C7, D7).
$ws.Range("D7").Value = "W"
$ws.Range("G7").Value = "-"

# --- Row 8: D8 holds the severity fed into VALIDATE(ISERROR(#VALUE!), B8, C8, D8).
$ws.Range("D8").Value = "E"
$ws.Range("G8").Value = "-"

# --- Row 9: D9 holds the severity fed into VALIDATE(ISLOGICAL(TRUE), B9, C9, D9).
$ws.Range("D9").Value = "W"
$ws.Range("G9").Value = "-"

# --- Row 10: D10 holds the severity fed into VALIDATE(ISERROR(E10), B10, C10, D10).
$ws.Range("D10").Value = "E"
$ws.Range("G10").Value = "Message A10, E, C10, #VALUE!, #VALUE!"

# --- Selection moved from C13 to G11.
$ws.Range("G11").Select()

# --- The workbook window was scrolled/repositioned (yWindow 555 -> 1455).
$excel.ActiveWindow.Top = 1455
